# Update of December 16th 21:00hr 10 minute Raw series stdDev values
# (wrong column had been pasted into "Bitcoin Raw Series StdDev" / "Bitcoin
# Returns Series StdDev" for the 16:00-21:00 10-minute-level rows).

$wb = $excel.ActiveWorkbook

$wsTenMin = $wb.Worksheets.Item("10 Minute Level")

# Column E = "Bitcoin Raw Series StdDev", Column G = "Bitcoin Returns Series StdDev"
# for rows 16-21 (16:00 - 21:00). The correct Raw-series StdDev values were
# pasted into the wrong cells previously - restore both columns to the
# right values.
$wsTenMin.Range("E16").Value = 43.889699999999998
$wsTenMin.Range("G16").Value = 0.0669

$wsTenMin.Range("E17").Value = 5.9950000000000001
$wsTenMin.Range("G17").Value = 0.1155

$wsTenMin.Range("E18").Value = 13.314399999999999
$wsTenMin.Range("G18").Value = 0.0929

$wsTenMin.Range("E19").Value = 12.843400000000001
$wsTenMin.Range("G19").Value = 0.1028

$wsTenMin.Range("E20").Value = 13.039099999999999
$wsTenMin.Range("G20").Value = 0.108

$wsTenMin.Range("E21").Value = 15.930400000000001
$wsTenMin.Range("G21").Value = 0.0744

# Update the active selection/tab: move the active tab from
# "10 Minute Level" back to "Day Level", and update each sheet's
# remembered selection.
[void]$wsTenMin.Range("H21").Select()

$wsDay = $wb.Worksheets.Item("Day Level")
$wsDay.Activate()
[void]$wsDay.Range("A7").Select()

# Reposition/resize the workbook window to match the saved view state.
$win = $wb.Windows.Item(1)
$win.Left = 2730
$win.Top = 2730
$win.Width = 28800
$win.Height = 18345
